$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style swap: header labels (B4:D4) gain left/top alignment; the
#     "No."/"Title" data columns (and the unused "Date" placeholders in
#     rows 24-30) lose their left/top alignment, keeping only wrap text.
$r4 = $ws.Range("B4:D4")
$r4.HorizontalAlignment = -4131; $r4.VerticalAlignment = -4160   # xlLeft / xlTop

$rb = $ws.Range("B5:B30")
$rb.HorizontalAlignment = 1; $rb.VerticalAlignment = -4107        # xlGeneral / xlBottom

$rd = $ws.Range("D5:D30")
$rd.HorizontalAlignment = 1; $rd.VerticalAlignment = -4107

$rc = $ws.Range("C24:C30")
$rc.HorizontalAlignment = 1; $rc.VerticalAlignment = -4107

# --- Correct the year typo on entry #6 (06/01/2020 -> 06/01/2021)
$ws.Range("C10").Value = 44202

# --- New description for entry #6
$ws.Range("E10").Value = "Just a rectangle to begin with."

# --- New logbook entry #7 (18/01/2021)
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 44214
$ws.Range("D11").Value = "Fixed bug, began implementing player movement"
$ws.Range("E11").Value = "Rectangle didn't display upon loading. Fixed this error and implemented player movement using the WASD keys."
$ws.Rows.Item(11).RowHeight = 29

# --- Move the selection to the newly-added row
$ws.Range("E11:G11").Select() | Out-Null
